$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 174
}

$ws.Range("C2:C16").Select()
